$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

$ws.Range("E9").Value = "3x0"
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = "Finalizado"

$ws.Range("E10").Value = "8x1"
$ws.Range("G10").Value = 8
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "Finalizado"

$ws.Columns.Item(9).ColumnWidth = 9.14

$ws.Range("L9").Select()
